# Insert a new data row at row 86 (pushing the existing row 86 and all
# rows below it down by one), then populate the new row with the
# observation that was added to the "Frambuesa" price series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 86..155 down to 87..156, leaving a blank row at 86.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86.
$ws.Range("A86").Value = 6
$ws.Range("B86").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C86").Value = "Metropolitana"
$ws.Range("D86").Value = 44589
$ws.Range("E86").Value = 13
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100101
$ws.Range("H86").Value = "Berries"
$ws.Range("I86").Value = 100101004
$ws.Range("J86").Value = "Frambuesa"
$ws.Range("K86").Value = "Sin especificar"
$ws.Range("L86").Value = "Especial"
$ws.Range("M86").Value = 500
$ws.Range("N86").Value = 8000
$ws.Range("O86").Value = 8000
$ws.Range("P86").Value = 8000
$ws.Range("Q86").Value = '$/bandeja 2 kilos'
$ws.Range("R86").Value = "Provincia de Linares"
$ws.Range("S86").Value = 4000
$ws.Range("T86").Value = 2
